$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.693.86"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.097.97"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5191"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4377"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09218"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.166"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.785"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "2.056.27"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.159"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "103.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.28%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06669"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.212"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "29.716.09"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "2.313.44"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.491"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.127"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.695"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.192"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.949"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.349"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06719"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.327"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6785"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.327"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000361"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.619"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.200"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
